$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clone row 8's formatting into the four new rows (9-12) first, while D8
#     still carries its original "no-fill" style (so the copy lands on D12,
#     the one new row that stays unfilled, without extra style churn). ---
$ws.Range("A8:E8").Copy()
$ws.Range("A9:E12").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Row 8: "Alert handling Validation" row - Execute (D8) flips from "Yes" to "No" ---
$ws.Cells.Item(8, 4).Value2 = "No"
$ws.Cells.Item(8, 4).Interior.Color = 5296274

# Row 9: Web Table Handling_Get coursename based on price
$ws.Cells.Item(9, 1).Value2 = 8
$ws.Cells.Item(9, 2).Value2 = "Web Table Handling_Get coursename based on price"
$ws.Cells.Item(9, 3).Value2 = "features/UnitTestCases.feature:31"
$ws.Cells.Item(9, 4).Value2 = "No"
$ws.Cells.Item(9, 4).Interior.Color = 5296274
$ws.Cells.Item(9, 5).Value2 = "PASSED"

# Row 10: Fixed Web Table Handling_Validate total price
$ws.Cells.Item(10, 1).Value2 = 9
$ws.Cells.Item(10, 2).Value2 = "Fixed Web Table Handling_Validate total price"
$ws.Cells.Item(10, 3).Value2 = "features/UnitTestCases.feature:35"
$ws.Cells.Item(10, 4).Value2 = "No"
$ws.Cells.Item(10, 4).Interior.Color = 5296274
$ws.Cells.Item(10, 5).Value2 = "PASSED"

# Row 11: Mouse Hover validation
$ws.Cells.Item(11, 1).Value2 = 10
$ws.Cells.Item(11, 2).Value2 = "Mouse Hover validation"
$ws.Cells.Item(11, 3).Value2 = "features/UnitTestCases.feature:39"
$ws.Cells.Item(11, 4).Value2 = "No"
$ws.Cells.Item(11, 4).Interior.Color = 5296274
$ws.Cells.Item(11, 5).Value2 = "PASSED"

# Row 12: Frame Validation (Execute stays "Yes", unfilled like the old row 8 was)
$ws.Cells.Item(12, 1).Value2 = 11
$ws.Cells.Item(12, 2).Value2 = "Frame Validation"
$ws.Cells.Item(12, 3).Value2 = "features/UnitTestCases.feature:43"
$ws.Cells.Item(12, 4).Value2 = "Yes"
$ws.Cells.Item(12, 5).Value2 = "PASSED"

# --- Widen column B to fit the new, longer test-case names ---
$ws.Columns(2).ColumnWidth = 54.57

# --- Selection ends on B12, matching the final edit position ---
$ws.Range("B12").Select()
